$wb = $excel.ActiveWorkbook

# OFF sheet: row 2 values updated (Week 15 simulations added)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 202
$wsOff.Range("C2").Value = 128
$wsOff.Range("D2").Value = 43
$wsOff.Range("E2").Value = 17
$wsOff.Range("F2").Value = 6
$wsOff.Range("G2").Value = 4

# DEF sheet: row 2 values updated (Week 15 simulations added)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 94
$wsDef.Range("C2").Value = 69
$wsDef.Range("D2").Value = 23
$wsDef.Range("E2").Value = 11
